# sample-jxls7.xlsx: update the jxls placeholder tokens to the lower-cased
# bean-property names used by the newer jxls7 template, and move the
# active selection the way the author last left it in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# jxls token ITEM_NAME -> item_name
$ws.Range("A5").Value = '${row.item_name}'
# jxls token VOLUME -> volume
$ws.Range("B5").Value = '${row.volume}'

# Selection moved from A5 to B6
[void]$ws.Range("B6").Select()
